$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.109.03'
$ws.Range('E2').Value = '  -2.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.740.18'
$ws.Range('E3').Value = '  -5.47%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '508.63'
$ws.Range('E5').Value = '  -3.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.11'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.752.03'
$ws.Range('E9').Value = '  -5.36%  '
$ws.Range('E10').Value = '  +4.27%  '
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.222.88'
$ws.Range('E14').Value = '  -5.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.027.66'
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.92'
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.742.11'
$ws.Range('E18').Value = '  -5.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.77'
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.09'
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.71'
$ws.Range('E21').Value = '  -3.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.30'
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.34'
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.175'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.428'
$ws.Range('E27').Value = '  -3.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0842'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.56'
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.27'
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.03'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.23'
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.46'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.963'
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.15'
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.23'
$ws.Range('E39').Value = '  -4.30%  '
$ws.Range('E40').Value = '  -4.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.56'
$ws.Range('E41').Value = '  -2.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.188.48'
$ws.Range('E42').Value = '  -6.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0561'
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.607'
$ws.Range('E45').Value = '  -5.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.22'
$ws.Range('E46').Value = '  -6.39%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0890'
$ws.Range('E50').Value = '  -3.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.28'
$ws.Range('E51').Value = '  +0.99%  '
